# Insert a new data row above row 384 (pushing the existing rows 384-408
# down to 385-409) and populate it with the new weekly price observation
# for "Paine" zapallo (1a nueva(o)).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(384).Insert()

$ws.Cells.Item(384, 1).Value = 4
$ws.Cells.Item(384, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(384, 3).Value = 'Los Lagos'
$ws.Cells.Item(384, 4).Value = 44931
$ws.Cells.Item(384, 5).Value = 10
$ws.Cells.Item(384, 6).Value = 100112045
$ws.Cells.Item(384, 7).Value = 'Zapallo'
$ws.Cells.Item(384, 8).Value = 'Paine'
$ws.Cells.Item(384, 9).Value = '1a nueva(o)'
$ws.Cells.Item(384, 10).Value = 500
$ws.Cells.Item(384, 11).Value = 650
$ws.Cells.Item(384, 12).Value = 700
$ws.Cells.Item(384, 13).Value = 675
$ws.Cells.Item(384, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(384, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(384, 16).Value = 675
$ws.Cells.Item(384, 17).Value = 1
$ws.Cells.Item(384, 18).Value = 'Hortaliza'
